$d = $word.ActiveDocument

# --- 1. Insert new sentence about the "sum_cols_first" option right after
#        the existing "If the option" lead-in (before "cond" is introduced). ---
$findRange = $d.Content
$findRange.Find.Execute(
    "If the option",  # FindText
    $false,            # MatchCase
    $false,            # MatchWholeWord
    $false,            # MatchWildcards
    $false,            # MatchSoundsLike
    $false,            # MatchAllWordForms
    $true,             # Forward
    1,                 # Wrap (wdFindContinue)
    $false,            # Format
    "",                # ReplaceWith (unused, we just locate the text)
    0                  # Replace (wdReplaceNone)
) | Out-Null

$findRange.Collapse(0)
$insertStart = $findRange.Start
$findRange.InsertAfter(" sum_cols_first is given, columns with summary statistics appear before columns with denominators. If the option")
$insertEnd = $findRange.End

# Re-apply the "VerbatimChar" character style to just the option name we
# inserted (scope the Find to the freshly-inserted text only, so the
# existing "sum_cols_first" occurrence in the heading is left untouched).
$scoped = $d.Range($insertStart, $insertEnd)
$scoped.Find.Execute("sum_cols_first", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$scoped.Style = "VerbatimChar"

# --- 2. Update the posted output line: second group now labelled "0"
#        instead of "2" for both the summary-statistic and N columns.
#        (Search/replace text deliberately excludes the surrounding
#        straight quotes so AutoFormat can't smart-quote them.) ---
$d.Content.Find.Execute("Summary 2", $false, $false, $false, $false, $false, $true, 1, $false, "Summary 0", 2) | Out-Null
$d.Content.Find.Execute("N 2", $false, $false, $false, $false, $false, $true, 1, $false, "N 0", 2) | Out-Null
